# Auto-generated edit script: apply numeric value updates to match target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I2").Value = 206.8
$ws.Range("H2").Value = 319
$ws.Range("M2").Value = -93.80000000000001
$ws.Range("K2").Value = 206.8
$ws.Range("K12").Value = 1155.6
$ws.Range("I12").Value = 1155.6
$ws.Range("H12").Value = 1098.625
$ws.Range("M12").Value = -985.5999999999999
$ws.Range("J51").Value = 4207.4165
$ws.Range("M51").Value = -2880.5
$ws.Range("L51").Value = 4207.4165
$ws.Range("N51").Value = -5175.4165
$ws.Range("K51").Value = 3364.5
$ws.Range("I51").Value = 3364.5
$ws.Range("H51").Value = 3824.2727
$ws.Range("L80").Value = 12023.1819
$ws.Range("I80").Value = 1547.0555
$ws.Range("K80").Value = 4641.166499999999
$ws.Range("J80").Value = 4007.7273
$ws.Range("H80").Value = 2480.4138
$ws.Range("N80").Value = -14019.1819
$ws.Range("M80").Value = -3643.166499999999
$ws.Range("I83").Value = 1547.0555
$ws.Range("N83").Value = -46053.5457
$ws.Range("H83").Value = 2480.4138
$ws.Range("M83").Value = -8931.4995
$ws.Range("K83").Value = 13923.4995
$ws.Range("J83").Value = 4007.7273
$ws.Range("L83").Value = 36069.5457
$ws.Range("M86").Value = -23809551
$ws.Range("I86").Value = 23810674
$ws.Range("L86").Value = 142859660
$ws.Range("J86").Value = 142859660
$ws.Range("H86").Value = 63493668
$ws.Range("K86").Value = 23810674
$ws.Range("N86").Value = -142861906
$ws.Range("I89").Value = 23810674
$ws.Range("N89").Value = -714309532
$ws.Range("J89").Value = 142859660
$ws.Range("K89").Value = 119053370
$ws.Range("H89").Value = 63493668
$ws.Range("M89").Value = -119047754
$ws.Range("L89").Value = 714298300
$ws.Range("H92").Value = 315.43478
$ws.Range("I92").Value = 322.75
$ws.Range("L92").Value = 266.66666
$ws.Range("M92").Value = 925.25
$ws.Range("K92").Value = 322.75
$ws.Range("N92").Value = -2762.66666
$ws.Range("J92").Value = 266.66666
$ws.Range("L96").Value = 1197
$ws.Range("N96").Value = -3943
$ws.Range("H96").Value = 91394.91
$ws.Range("J96").Value = 399
$ws.Range("K98").Value = 934.2857
$ws.Range("I98").Value = 934.2857
$ws.Range("H98").Value = 934.2857
$ws.Range("M98").Value = 563.7143
$ws.Range("M99").Value = -3737
$ws.Range("H99").Value = 1889.2858
$ws.Range("I99").Value = 1745
$ws.Range("L99").Value = 5740.0002
$ws.Range("K99").Value = 5235
$ws.Range("J99").Value = 1913.3334
$ws.Range("N99").Value = -8736.0002
$ws.Range("I100").Value = 5000
$ws.Range("L100").Value = 2484.1667
$ws.Range("M100").Value = -4459
$ws.Range("K100").Value = 5000
$ws.Range("N100").Value = -3566.1667
$ws.Range("J100").Value = 2484.1667
$ws.Range("H100").Value = 2843.5715
$ws.Range("N108").Value = -92671.664
$ws.Range("L108").Value = 84991.664
$ws.Range("J108").Value = 84991.664
$ws.Range("H108").Value = 84991.664
$ws.Range("J114").Value = 99941.164
$ws.Range("L114").Value = 99941.164
$ws.Range("H114").Value = 99941.164
$ws.Range("N114").Value = -108619.164
$ws.Range("N116").Value = -3046709.8
$ws.Range("I116").Value = 8785.286
$ws.Range("J116").Value = 3039825.8
$ws.Range("L116").Value = 3039825.8
$ws.Range("K116").Value = 8785.286
$ws.Range("H116").Value = 1861087.8
$ws.Range("M116").Value = -5343.286
$ws.Range("H117").Value = 88082.75
$ws.Range("J117").Value = 88082.75
$ws.Range("N117").Value = -97260.75
$ws.Range("L117").Value = 88082.75
$ws.Range("M118").Value = -199.1428999999998
$ws.Range("I118").Value = 618.7143
$ws.Range("K118").Value = 1856.1429
$ws.Range("H118").Value = 618.7143
$ws.Range("J120").Value = 38128.168
$ws.Range("N120").Value = -47804.168
$ws.Range("H120").Value = 38128.168
$ws.Range("L120").Value = 38128.168
$ws.Range("H122").Value = 934.2857
$ws.Range("I122").Value = 934.2857
$ws.Range("M122").Value = -352.8571000000002
$ws.Range("K122").Value = 2802.8571
$ws.Range("J123").Value = 89448
$ws.Range("H123").Value = 87873.336
$ws.Range("L123").Value = 89448
$ws.Range("N123").Value = -99248
$ws.Range("L125").Value = 60260.724
$ws.Range("H125").Value = 5581.3887
$ws.Range("N125").Value = -65180.724
$ws.Range("K125").Value = 34473.8565
$ws.Range("M125").Value = -32013.8565
$ws.Range("I125").Value = 3830.4285
$ws.Range("J125").Value = 6695.636
$ws.Range("I127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("K127").Value = 0
$ws.Range("H127").Value = 0
$ws.Range("H131").Value = 1823
$ws.Range("I131").Value = 1823
$ws.Range("M131").Value = -429
$ws.Range("K131").Value = 5469
$ws.Range("K132").Value = 3731.298
$ws.Range("N132").Value = -13190
$ws.Range("J132").Value = 2710
$ws.Range("H132").Value = 1521.8448
$ws.Range("L132").Value = 8130
$ws.Range("M132").Value = -1201.298
$ws.Range("I132").Value = 1243.766
$ws.Range("J134").Value = 57913.332
$ws.Range("H134").Value = 57913.332
$ws.Range("N134").Value = -68053.33199999999
$ws.Range("L134").Value = 57913.332
$ws.Range("N136").Value = -89071.336
$ws.Range("L136").Value = 78871.336
$ws.Range("H136").Value = 78871.336
$ws.Range("J136").Value = 78871.336
$ws.Range("I137").Value = 1948.7222
$ws.Range("H137").Value = 581936.1
$ws.Range("K137").Value = 5846.1666
$ws.Range("M137").Value = -3296.1666
$ws.Range("N139").Value = -110265
$ws.Range("L139").Value = 99985
$ws.Range("H139").Value = 99985
$ws.Range("J139").Value = 99985

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J16").Value = 50
$ws.Range("I16").Value = 697.5
$ws.Range("L16").Value = 50
$ws.Range("K16").Value = 697.5
$ws.Range("N16").Value = -624
$ws.Range("H16").Value = 481.66666
$ws.Range("M16").Value = -410.5
$ws.Range("J31").Value = 65333
$ws.Range("N31").Value = -65921
$ws.Range("H31").Value = 33183
$ws.Range("L31").Value = 65333
$ws.Range("I32").Value = 6332.706
$ws.Range("K32").Value = 6332.706
$ws.Range("L32").Value = 30224.084
$ws.Range("M32").Value = -6045.706
$ws.Range("J32").Value = 30224.084
$ws.Range("H32").Value = 16218.793
$ws.Range("N32").Value = -30798.084
$ws.Range("J52").Value = 54364
$ws.Range("L52").Value = 54364
$ws.Range("N52").Value = -55000
$ws.Range("H52").Value = 54364
$ws.Range("H61").Value = 26861.025
$ws.Range("L61").Value = 115791.78
$ws.Range("M61").Value = -830.4193
$ws.Range("K61").Value = 1042.4193
$ws.Range("I61").Value = 1042.4193
$ws.Range("J61").Value = 115791.78
$ws.Range("N61").Value = -116215.78
$ws.Range("M74").Value = -77811.92
$ws.Range("J74").Value = 2792.5334
$ws.Range("K74").Value = 78685.92
$ws.Range("L74").Value = 2792.5334
$ws.Range("H74").Value = 38028.75
$ws.Range("I74").Value = 78685.92
$ws.Range("N74").Value = -4540.5334
$ws.Range("J77").Value = 2792.5334
$ws.Range("I77").Value = 78685.92
$ws.Range("L77").Value = 13962.667
$ws.Range("H77").Value = 38028.75
$ws.Range("N77").Value = -22698.667
$ws.Range("M77").Value = -389061.6
$ws.Range("K77").Value = 393429.6
$ws.Range("H97").Value = 1656
$ws.Range("M97").Value = -948.875
$ws.Range("I97").Value = 1444.875
$ws.Range("K97").Value = 1444.875
$ws.Range("I102").Value = 35745.656
$ws.Range("L102").Value = 27377.5
$ws.Range("N102").Value = -30621.5
$ws.Range("M102").Value = -34123.656
$ws.Range("J102").Value = 27377.5
$ws.Range("K102").Value = 35745.656
$ws.Range("H102").Value = 34731.332
$ws.Range("K110").Value = 1214.5714
$ws.Range("M110").Value = 830.4286
$ws.Range("I110").Value = 1214.5714
$ws.Range("H110").Value = 1175.25
$ws.Range("L115").Value = 53200.668
$ws.Range("N115").Value = -56334.668
$ws.Range("J115").Value = 53200.668
$ws.Range("H115").Value = 49932.145
$ws.Range("L119").Value = 0
$ws.Range("H119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("J119").Value = 0
$ws.Range("K132").Value = 3499.8
$ws.Range("N132").Value = -14058.3329
$ws.Range("J132").Value = 2999.4443
$ws.Range("H132").Value = 1541.5
$ws.Range("L132").Value = 8998.332900000001
$ws.Range("M132").Value = -969.7999999999997
$ws.Range("I132").Value = 1166.6
$ws.Range("N136").Value = -352475.34
$ws.Range("L136").Value = 347375.34
$ws.Range("H136").Value = 26861.025
$ws.Range("I136").Value = 1042.4193
$ws.Range("J136").Value = 115791.78
$ws.Range("K136").Value = 3127.2579
$ws.Range("M136").Value = -577.2579000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L2").Value = 28081.111
$ws.Range("N2").Value = -28307.111
$ws.Range("H2").Value = 28081.111
$ws.Range("J2").Value = 28081.111
$ws.Range("K20").Value = 112846.336
$ws.Range("H20").Value = 61092.35
$ws.Range("M20").Value = -112599.336
$ws.Range("I20").Value = 112846.336
$ws.Range("J51").Value = 43993.332
$ws.Range("L51").Value = 43993.332
$ws.Range("N51").Value = -44975.332
$ws.Range("H51").Value = 43993.332
$ws.Range("H94").Value = 2702.75
$ws.Range("I94").Value = 3665.8462
$ws.Range("K94").Value = 3665.8462
$ws.Range("M94").Value = -3214.8462
$ws.Range("M99").Value = -60140.707
$ws.Range("H99").Value = 1468895.9
$ws.Range("I99").Value = 61638.707
$ws.Range("L99").Value = 6253570
$ws.Range("K99").Value = 61638.707
$ws.Range("J99").Value = 6253570
$ws.Range("N99").Value = -6256566
$ws.Range("J105").Value = 5625
$ws.Range("N105").Value = -9119
$ws.Range("K105").Value = 8501583
$ws.Range("L105").Value = 5625
$ws.Range("H105").Value = 5103200
$ws.Range("M105").Value = -8499836
$ws.Range("I105").Value = 8501583
$ws.Range("L110").Value = 99990
$ws.Range("J110").Value = 99990
$ws.Range("H110").Value = 99990
$ws.Range("N110").Value = -108170
$ws.Range("L115").Value = 82115.75
$ws.Range("N115").Value = -85249.75
$ws.Range("J115").Value = 82115.75
$ws.Range("H115").Value = 79880.78
$ws.Range("L118").Value = 45082.5
$ws.Range("N118").Value = -48396.5
$ws.Range("H118").Value = 46245.11
$ws.Range("J118").Value = 45082.5
$ws.Range("L119").Value = 42652.25
$ws.Range("H119").Value = 42652.25
$ws.Range("N119").Value = -52328.25
$ws.Range("J119").Value = 42652.25
$ws.Range("L127").Value = 89944.336
$ws.Range("H127").Value = 89944.336
$ws.Range("J127").Value = 89944.336
$ws.Range("N127").Value = -99864.336
$ws.Range("N132").Value = -91462.86
$ws.Range("J132").Value = 81342.86
$ws.Range("H132").Value = 81342.86
$ws.Range("L132").Value = 81342.86
$ws.Range("K134").Value = 3334.6155
$ws.Range("I134").Value = 1111.5385
$ws.Range("H134").Value = 1739.1702
$ws.Range("M134").Value = -799.6155000000003
$ws.Range("H135").Value = 81764
$ws.Range("L135").Value = 81764
$ws.Range("J135").Value = 81764
$ws.Range("N135").Value = -91904
$ws.Range("J138").Value = 88003.39999999999
$ws.Range("H138").Value = 88003.39999999999
$ws.Range("N138").Value = -98283.39999999999
$ws.Range("L138").Value = 88003.39999999999
$ws.Range("J140").Value = 65911.75
$ws.Range("N140").Value = -76271.75
$ws.Range("H140").Value = 114143.78
$ws.Range("L140").Value = 65911.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M31").Value = -1687.0714
$ws.Range("I31").Value = 1982.0714
$ws.Range("H31").Value = 3158.3845
$ws.Range("K31").Value = 1982.0714
$ws.Range("M34").Value = -1780.0714
$ws.Range("K34").Value = 1982.0714
$ws.Range("H34").Value = 3158.3845
$ws.Range("I34").Value = 1982.0714
$ws.Range("K105").Value = 2494.4
$ws.Range("H105").Value = 4241.222
$ws.Range("M105").Value = -747.4000000000001
$ws.Range("I105").Value = 2494.4
$ws.Range("N116").Value = -59154.125
$ws.Range("J116").Value = 49976.125
$ws.Range("L116").Value = 49976.125
$ws.Range("H116").Value = 49976.125
$ws.Range("H117").Value = 42196.445
$ws.Range("J117").Value = 42196.445
$ws.Range("N117").Value = -51374.445
$ws.Range("L117").Value = 42196.445
$ws.Range("H122").Value = 2198.975
$ws.Range("I122").Value = 1910.0769
$ws.Range("M122").Value = -3280.2307
$ws.Range("K122").Value = 5730.2307
$ws.Range("K132").Value = 4567.9092
$ws.Range("H132").Value = 1788.3846
$ws.Range("M132").Value = -2037.9092
$ws.Range("I132").Value = 1522.6364
$ws.Range("K134").Value = 5913.6
$ws.Range("I134").Value = 1971.2
$ws.Range("H134").Value = 42078
$ws.Range("M134").Value = -3378.6
$ws.Range("J138").Value = 99840
$ws.Range("H138").Value = 99840
$ws.Range("N138").Value = -110120
$ws.Range("L138").Value = 99840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M40").Value = -135.615388
$ws.Range("I40").Value = 51.153847
$ws.Range("J40").Value = 147.8
$ws.Range("K40").Value = 204.615388
$ws.Range("H40").Value = 78
$ws.Range("N40").Value = -729.2
$ws.Range("L40").Value = 591.2
$ws.Range("K109").Value = 870
$ws.Range("J109").Value = 4000
$ws.Range("I109").Value = 290
$ws.Range("N109").Value = -14080
$ws.Range("H109").Value = 2145
$ws.Range("M109").Value = 170
$ws.Range("L109").Value = 12000
$ws.Range("K110").Value = 12375
$ws.Range("M110").Value = -8285
$ws.Range("I110").Value = 4125
$ws.Range("H110").Value = 5475
$ws.Range("M112").Value = -9086
$ws.Range("H112").Value = 5127.75
$ws.Range("I112").Value = 3398
$ws.Range("K112").Value = 10194
$ws.Range("N116").Value = -10362.5
$ws.Range("J116").Value = 1159.5
$ws.Range("L116").Value = 3478.5
$ws.Range("H116").Value = 1100.7778
$ws.Range("M119").Value = 2528.75
$ws.Range("H119").Value = 3615.8
$ws.Range("K119").Value = 2309.25
$ws.Range("I119").Value = 769.75
$ws.Range("J120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H120").Value = 15185.6
$ws.Range("L120").Value = 0
$ws.Range("H122").Value = 2020718
$ws.Range("N122").Value = -18191362
$ws.Range("L122").Value = 18186462
$ws.Range("J122").Value = 2020718
$ws.Range("I139").Value = 1671.8422
$ws.Range("H139").Value = 2464.9524
$ws.Range("M139").Value = 124.4733999999999
$ws.Range("K139").Value = 5015.5266

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 5289.8
$ws.Range("J22").Value = 6483.3335
$ws.Range("N22").Value = -7541.3335
$ws.Range("K22").Value = 3499.5
$ws.Range("I22").Value = 3499.5
$ws.Range("M22").Value = -2970.5
$ws.Range("L22").Value = 6483.3335
$ws.Range("L43").Value = 19666.334
$ws.Range("N43").Value = -19968.334
$ws.Range("H43").Value = 12813.125
$ws.Range("J43").Value = 19666.334
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 335169.34
$ws.Range("M70").Value = -334899.34
$ws.Range("H70").Value = 252127
$ws.Range("L70").Value = 3000
$ws.Range("I70").Value = 335169.34
$ws.Range("N70").Value = -3540
$ws.Range("J73").Value = 3000
$ws.Range("L73").Value = 3000
$ws.Range("H73").Value = 252127
$ws.Range("N73").Value = -4872
$ws.Range("K73").Value = 335169.34
$ws.Range("I73").Value = 335169.34
$ws.Range("M73").Value = -334233.34
$ws.Range("H97").Value = 2669.7646
$ws.Range("L97").Value = 10999.667
$ws.Range("M97").Value = -388.7857
$ws.Range("N97").Value = -11991.667
$ws.Range("J97").Value = 10999.667
$ws.Range("I97").Value = 884.7857
$ws.Range("K97").Value = 884.7857
$ws.Range("I102").Value = 1502.5
$ws.Range("M102").Value = 119.5
$ws.Range("K102").Value = 1502.5
$ws.Range("H102").Value = 1502.5
$ws.Range("K108").Value = 60000
$ws.Range("N108").Value = -69968.25
$ws.Range("L108").Value = 62288.25
$ws.Range("M108").Value = -56160
$ws.Range("I108").Value = 60000
$ws.Range("J108").Value = 62288.25
$ws.Range("H108").Value = 62112.23
$ws.Range("L119").Value = 67567.73
$ws.Range("H119").Value = 67567.73
$ws.Range("N119").Value = -77243.73
$ws.Range("J119").Value = 67567.73
$ws.Range("K132").Value = 6327.5625
$ws.Range("H132").Value = 2926.5356
$ws.Range("M132").Value = -3797.5625
$ws.Range("I132").Value = 2109.1875
$ws.Range("H135").Value = 94277.14
$ws.Range("L135").Value = 94277.14
$ws.Range("J135").Value = 94277.14
$ws.Range("N135").Value = -104417.14
$ws.Range("J140").Value = 57272.25
$ws.Range("N140").Value = -67632.25
$ws.Range("H140").Value = 55017.8
$ws.Range("L140").Value = 57272.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 687.2917
$ws.Range("J22").Value = 769.6875
$ws.Range("N22").Value = -1359.6875
$ws.Range("L22").Value = 769.6875
$ws.Range("L27").Value = 769.6875
$ws.Range("N27").Value = -983.6875
$ws.Range("H27").Value = 687.2917
$ws.Range("J27").Value = 769.6875
$ws.Range("L82").Value = 1939.75
$ws.Range("I82").Value = 1423.5
$ws.Range("H82").Value = 1571
$ws.Range("K82").Value = 1423.5
$ws.Range("J82").Value = 1939.75
$ws.Range("N82").Value = -2661.75
$ws.Range("M82").Value = -1062.5
$ws.Range("L85").Value = 1939.75
$ws.Range("J85").Value = 1939.75
$ws.Range("N85").Value = -4435.75
$ws.Range("M85").Value = -175.5
$ws.Range("H85").Value = 1571
$ws.Range("K85").Value = 1423.5
$ws.Range("I85").Value = 1423.5
$ws.Range("H93").Value = 1735.8889
$ws.Range("L93").Value = 2149.1667
$ws.Range("J93").Value = 2149.1667
$ws.Range("N93").Value = -4645.1667
$ws.Range("J121").Value = 52978.89
$ws.Range("N121").Value = -56472.89
$ws.Range("L121").Value = 52978.89
$ws.Range("H121").Value = 52978.89
$ws.Range("J123").Value = 79991.11
$ws.Range("H123").Value = 79991.11
$ws.Range("L123").Value = 79991.11
$ws.Range("N123").Value = -89791.11
$ws.Range("K132").Value = 27225
$ws.Range("N132").Value = -12529.25
$ws.Range("J132").Value = 2489.75
$ws.Range("H132").Value = 7428.6875
$ws.Range("L132").Value = 7469.25
$ws.Range("M132").Value = -24695
$ws.Range("I132").Value = 9075
$ws.Range("N136").Value = -15470.88
$ws.Range("L136").Value = 10370.88
$ws.Range("H136").Value = 3116.558
$ws.Range("I136").Value = 2643.7778
$ws.Range("J136").Value = 3456.96
$ws.Range("K136").Value = 7931.3334
$ws.Range("M136").Value = -5381.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L96").Value = 3125
$ws.Range("N96").Value = -5871
$ws.Range("H96").Value = 2741.6667
$ws.Range("J96").Value = 3125
$ws.Range("H107").Value = 1590.3611
$ws.Range("L107").Value = 6879.176399999999
$ws.Range("M107").Value = -964.8948
$ws.Range("J107").Value = 2293.0588
$ws.Range("I107").Value = 961.6316
$ws.Range("K107").Value = 2884.8948
$ws.Range("N107").Value = -10719.1764
$ws.Range("J121").Value = 49073.332
$ws.Range("N121").Value = -52567.332
$ws.Range("L121").Value = 49073.332
$ws.Range("H121").Value = 49073.332
$ws.Range("H122").Value = 3537.1667
$ws.Range("I122").Value = 3084
$ws.Range("M122").Value = -6802
$ws.Range("K122").Value = 9252
$ws.Range("M126").Value = -2670.821599999999
$ws.Range("I126").Value = 1713.6072
$ws.Range("H126").Value = 1854.2258
$ws.Range("K126").Value = 5140.821599999999
$ws.Range("N136").Value = -12768
$ws.Range("L136").Value = 7668
$ws.Range("H136").Value = 1716.6604
$ws.Range("I136").Value = 1544.9773
$ws.Range("J136").Value = 2556
$ws.Range("K136").Value = 4634.9319
$ws.Range("M136").Value = -2084.9319

